$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.523.85'
$ws.Range("E2").Value = '  +0.45%  '

$ws.Range("D3").Value = '2.466.43'
$ws.Range("E3").Value = '  -0.61%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.14%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '91.50'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.25%  '

$ws.Range("E7").Value = '  -0.87%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("E9").Value = '  +2.50%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.10'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.39%  '

$ws.Range("E11").Value = '  +0.46%  '

$ws.Range("E12").Value = '  +1.03%  '

$ws.Range("D13").Value = '2.846.14'
$ws.Range("E13").Value = '  -0.53%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.84'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.41%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.01'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.83%  '

$ws.Range("D16").Value = '2.493.26'
$ws.Range("E16").Value = '  +0.90%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.765'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.44%  '

$ws.Range("D18").Value = '41.487.61'
$ws.Range("E18").Value = '  +0.29%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.33%  '

$ws.Range("D20").Value = '0.0₃0947'
$ws.Range("E20").Value = '  +1.97%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.30'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.64%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.04'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.89%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.31'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.00%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.71'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.48%  '

$ws.Range("E25").Value = '  -0.11%  '

$ws.Range("E26").Value = '  -1.58%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.53'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.66%  '

$ws.Range("E28").Value = '  -0.52%  '

$ws.Range("E29").Value = '  -0.62%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.27'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '155.73'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.29%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.43'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.64%  '

$ws.Range("E33").Value = '  -0.55%  '

$ws.Range("E34").Value = '  +0.96%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.10'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.02%  '

$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.86'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.95%  '

$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.103'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.58%  '

$ws.Range("B38").Value = 'ApeXProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.29'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -10.57%  '

$ws.Range("E39").Value = '  -0.97%  '

$ws.Range("E40").Value = '  -4.92%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.35%  '

$ws.Range("D43").Value = '1.940.96'
$ws.Range("E43").Value = '  -2.73%  '

$ws.Range("E44").Value = '  -1.19%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.46'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.08%  '

$ws.Range("E46").Value = '  -3.47%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.02'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.10%  '

$ws.Range("D48").Value = '2.706.22'
$ws.Range("E48").Value = '  -0.66%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '96.88'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.55%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '66.71'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.25%  '

$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.58'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.01%  '
